$d = $word.ActiveDocument

# 1. "Proposed and led the first implementation of tracking user actions using
#    Datadog to better understand user needs and frustrations"
#    -> "...using Datadog to track KPIs and better understand user needs and frustrations"
$d.Content.Find.Execute(
    "Datadog to better understand",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "Datadog to track KPIs and better understand", 2) | Out-Null

# 2. "Co-owned the design implementations ... while rebranding the app"
#    -> "...while rebranding and redeveloping the app"
$d.Content.Find.Execute(
    "while rebranding the app",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "while rebranding and redeveloping the app", 2) | Out-Null

# 3. "Languages: ... Solidity, bash, SQL, HTML, CSS" -> capitalise "bash" to "Bash"
$d.Content.Find.Execute(
    "Solidity, bash, SQL",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "Solidity, Bash, SQL", 2) | Out-Null

# 4. "Technologies: ... Tailwind CSS, MongoDB, SQLite" -> swap order to "SQLite, MongoDB"
$d.Content.Find.Execute(
    "Tailwind CSS, MongoDB, SQLite",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "Tailwind CSS, SQLite, MongoDB", 2) | Out-Null

# 5. "DevOps: AWS, Git, Cloudflare Workers, GitHub Actions, Husky, Jest, Swagger"
#    -> "DevOps: AWS, Git, GitHub Actions, Datadog, Cloudflare Workers, Husky, Jest, Swagger"
$d.Content.Find.Execute(
    "Git, Cloudflare Workers, GitHub Actions, Husky, Jest, Swagger",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "Git, GitHub Actions, Datadog, Cloudflare Workers, Husky, Jest, Swagger", 2) | Out-Null
